$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.233599066734314
$ws.Range("B1").Value = 2.616581678390503
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.716055870056152
$ws.Range("E1").Value = 1.147258520126343
